$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.1072256642246516, 0.6857061106036189, 1.74092127272603, 3.7181689313671775)"
$ws.Range("C2").Value = "JSU(-1.2994241214309599, 1.238578647144636, 2.671559827800122, 4.733906329365339)"
$ws.Range("D2").Value = "NIG(0.6813670780981159, 0.506189197501962, 0.9207467703822829, 2.8778293875144274)"
$ws.Range("E2").Value = "JSU(-1.5385398959023293, 1.1378474706910915, 1.6161296508069531, 4.314194682566126)"
